$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.255104666666667
$ws.Range("N2").Value = 3.765314
$ws.Range("O2").Value = 0.08274843238311479
$ws.Range("P2").Value = 0.08274843238311479
$ws.Range("Q2").Value = 2.035564728067111
$ws.Range("R2").Value = 18.320082552604
$ws.Range("S2").Value = 0.08274843238311479
$ws.Range("T2").Value = 0.08274843238311479

# Row 3
$ws.Range("O3").Value = 0.1886083876486254
$ws.Range("P3").Value = 0.1886083876486254
$ws.Range("S3").Value = 0.1886083876486254
$ws.Range("T3").Value = 0.1886083876486254

# Row 4
$ws.Range("M4").Value = 8.752828666666668
$ws.Range("N4").Value = 26.258486
$ws.Range("O4").Value = 0.5770696821709866
$ws.Range("P4").Value = 0.5770696821709866
$ws.Range("Q4").Value = 14.19558844602178
$ws.Range("R4").Value = 127.760296014196
$ws.Range("S4").Value = 0.5770696821709866
$ws.Range("T4").Value = 0.5770696821709866

# Row 5
$ws.Range("M5").Value = 2.299023666666667
$ws.Range("N5").Value = 6.897071
$ws.Range("O5").Value = 0.1515734977972732
$ws.Range("P5").Value = 0.1515734977972732
$ws.Range("Q5").Value = 3.728622487945111
$ws.Range("R5").Value = 33.557602391506
$ws.Range("S5").Value = 0.1515734977972732
$ws.Range("T5").Value = 0.1515734977972732
